# Add a new "AnimalByproduct" column to the chemicals table (Table3).
# This records whether each chemical/reagent is an animal byproduct, which
# is sometimes required for regulatory purposes, and enables exporting the
# chemicals flagged as animal byproducts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Append a new column to the table (goes right after the last column,
# i.e. after "exitDate").
$newCol = $tbl.ListColumns.Add()

# Name the new column header.
$ws.Range("Q1").Value = "AnimalByproduct"

# Populate the data rows with default value "No".
$ws.Range("Q2").Value = "No"
$ws.Range("Q3").Value = "No"
$ws.Range("Q4").Value = "No"
